# Apply updated CO2 Change by Activities values (large-scale intervention run)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C2" = -0.00019764741600725169
    "D2" = -0.03895395257964651137
    "F2" = -0.00019764741600725169
    "G2" = -0.00098823708003648036
    "C3" = 0.00015789056953552369
    "D3" = 0.00063156227814209487
    "E3" = 0.00019304834286892891
    "F3" = 0.00029266487387169532
    "G3" = 0.00078945284767772961
    "C4" = 0.00001064393317176237
    "D4" = 0.00004257573268704948
    "E4" = 0.04868250540540942223
    "F4" = -13.26935995645544963395
    "G4" = 0.00005321966585881271
    "C5" = -0.44694350287136330735
    "D5" = -0.00075069064195121271
    "F5" = -0.0001876726604878032
    "G5" = -0.00093836330243846078
    "C6" = -0.00019171372591893221
    "D6" = -0.03778449290744489986
    "F6" = -0.00019171372591893221
    "G6" = -0.00095856862959464717
    "C7" = -0.00026358088734668161
    "D7" = -0.00105432354938672602
    "E7" = -1.20554852853280203284
    "F7" = -0.19495268395919398974
    "G7" = -0.00131790443673340807
    "B8" = -0.20867638562117460377
    "C8" = -0.00115358483246197395
    "D8" = -0.11246966312205590455
    "F8" = -0.00115358483246197395
    "G8" = -0.00576792416231342031
    "C9" = -0.00008342150408757742
    "D9" = -0.07188133868055501807
    "F9" = -0.00008342150408757742
    "G9" = -0.00041710752043755411
    "C10" = -0.00168261046852435903
    "D10" = -0.02566014622885859006
    "E10" = -0.0018337684911386991
    "F10" = -0.00376648224269615596
    "G10" = -0.00841305234263245438
    "C11" = -0.00032798194886396459
    "D11" = -0.0376787023833458079
    "F11" = -0.02089905984962570074
    "G11" = -0.2448995292497783971
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

